$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates: volume/number and reporting week dates
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 29   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/14/2022  Through  11/20/2022"

# ---------------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------------
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -66.666666666666
$ws.Range("N15").Value = 7.692307692307

# ---------------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 81
$ws.Range("J16").Value = 75
$ws.Range("K16").Value = 8
$ws.Range("L16").Value = 47.272727272727
$ws.Range("M16").Value = 44.642857142857
$ws.Range("N16").Value = -85.587188612099

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 14
$ws.Range("H17").Value = -12.5
$ws.Range("I17").Value = 105
$ws.Range("J17").Value = 102
$ws.Range("K17").Value = 2.941176470588
$ws.Range("L17").Value = 77.966101694915
$ws.Range("M17").Value = 110
$ws.Range("N17").Value = -26.573426573426

# ---------------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 66.666666666666
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -5.882352941176
$ws.Range("I18").Value = 165
$ws.Range("J18").Value = 136
$ws.Range("K18").Value = 21.323529411764
$ws.Range("L18").Value = -1.197604790419
$ws.Range("M18").Value = 96.428571428571
$ws.Range("N18").Value = -84.848484848484

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -18.75
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = 3.636363636363
$ws.Range("I19").Value = 591
$ws.Range("J19").Value = 484
$ws.Range("K19").Value = 22.107438016528
$ws.Range("L19").Value = 15.204678362573
$ws.Range("M19").Value = -10.725075528700
$ws.Range("N19").Value = -70.901033973412

# ---------------------------------------------------------------------------
# Row 20 -- C20/D20/E20 change data type (text <-> number), handled via a
# format+value paste from donor cells that already carry the desired style,
# so that the existing shared-string slots / numeric style ids are reused
# instead of creating brand-new styles or shared strings.
# ---------------------------------------------------------------------------
$ws.Range("F15").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value = 1

$ws.Range("C23").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163) | Out-Null
$ws.Range("C23").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null

$ws.Range("E23").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4163) | Out-Null
$ws.Range("E23").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null

$ws.Range("I20").Value = 59
$ws.Range("K20").Value = 34.090909090909
$ws.Range("L20").Value = 37.209302325581
$ws.Range("M20").Value = 96.666666666666
$ws.Range("N20").Value = -89.594356261022

# ---------------------------------------------------------------------------
# Row 21
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 100
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = -2.912621359223
$ws.Range("I21").Value = 1015
$ws.Range("J21").Value = 849
$ws.Range("K21").Value = 19.552414605418
$ws.Range("L21").Value = 19.693396226415
$ws.Range("M21").Value = 14.301801801801
$ws.Range("N21").Value = -76.958002270147

# ---------------------------------------------------------------------------
# Row 22 -- D22/E22 change data type (number -> text)
# ---------------------------------------------------------------------------
$ws.Range("C22").Value = 2

$ws.Range("D23").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("D23").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null

$ws.Range("E23").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4163) | Out-Null
$ws.Range("E23").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null

$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 28
$ws.Range("K22").Value = 27.272727272727
$ws.Range("L22").Value = 40
$ws.Range("M22").Value = -9.677419354838

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 57.142857142857
$ws.Range("F24").Value = 72
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = -2.702702702702
$ws.Range("I24").Value = 1130
$ws.Range("J24").Value = 884
$ws.Range("K24").Value = 27.828054298642
$ws.Range("L24").Value = -3.911564625850
$ws.Range("M24").Value = 101.426024955437

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
$ws.Range("F25").Value = 24
$ws.Range("H25").Value = 41.176470588235
$ws.Range("I25").Value = 213
$ws.Range("J25").Value = 159
$ws.Range("K25").Value = 33.962264150943
$ws.Range("L25").Value = 47.916666666666
$ws.Range("M25").Value = 2.898550724637

# ---------------------------------------------------------------------------
# Row 26 -- D26/E26 change data type (text -> number)
# ---------------------------------------------------------------------------
$ws.Range("D16").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").Value = 1

$ws.Range("E16").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = -100

$ws.Range("J26").Value = 13
$ws.Range("K26").Value = 23.076923076923

# ---------------------------------------------------------------------------
# Row 27
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 50
$ws.Range("J27").Value = 63
$ws.Range("K27").Value = -20.634920634920
$ws.Range("L27").Value = 25

# ---------------------------------------------------------------------------
# Row 30 -- G30/H30 change data type (number -> text)
# ---------------------------------------------------------------------------
$ws.Range("G23").Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4163) | Out-Null
$ws.Range("G23").Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4122) | Out-Null

$ws.Range("H23").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4163) | Out-Null
$ws.Range("H23").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
